# Apply the data-value updates to Sheet1 and move the active selection,
# matching the authored OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F ("avg days to close" or similar) value corrections.
$ws.Range("F7").Value = 822.79580111040002
$ws.Range("F10").Value = 718.58079578360002
$ws.Range("F25").Value = 487.54909416689998
$ws.Range("F27").Value = 567.49253346399996

# Move the active cell / selection from S16 to F27.
$ws.Activate()
$ws.Range("F27").Select()
